# Backlogs.xlsx — "Add files via upload": populate the previously-empty
# Sprint1 sheet with the "SPRINT BACKLOG 1" table, and restore the
# Product sheet's selection/active-cell state.

$wb = $excel.ActiveWorkbook
$product = $wb.Worksheets.Item("Product")
$ws = $wb.Worksheets.Item("Sprint1")

# ---------------------------------------------------------------------
# 1. Title row (copy the "PRODUCT BACKLOG" banner formatting so the new
#    title reuses the existing fill/font styles instead of minting new
#    ones) then overwrite the text.
# ---------------------------------------------------------------------
$product.Range("A1:F1").Copy()
$ws.Range("A1").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("A1").Value = "SPRINT BACKLOG 1"

# ---------------------------------------------------------------------
# 2. Column header row — reuse the existing header styles too.
# ---------------------------------------------------------------------
$product.Range("A2:F2").Copy()
$ws.Range("A2").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("A2").Value = "ID"
$ws.Range("B2").Value = "User story"
$ws.Range("C2").Value = "Priority"
$ws.Range("D2").Value = "Time estimation"
$ws.Range("E2").Value = "Tasks"
$ws.Range("F2").Value = "Acceptance criteria"

# ---------------------------------------------------------------------
# 3. Data rows (4 user stories pulled into sprint 1).
# ---------------------------------------------------------------------

# -- Row 3 : ID 0, menu with 4 options --------------------------------
$ws.Range("A3").Value = 0
$ws.Range("A3").HorizontalAlignment = -4108   # xlCenter
$ws.Range("A3").VerticalAlignment = -4108     # xlCenter

$ws.Range("B3").Value = "As a player I want to have a menu with 4 options, so `nthat I can control the game"
$ws.Range("B3").WrapText = $true
$ws.Range("B3").HorizontalAlignment = -4131   # xlLeft
$ws.Range("B3").VerticalAlignment = -4160     # xlTop

$ws.Range("C3").Value = "M"
$ws.Range("C3").HorizontalAlignment = -4108
$ws.Range("C3").VerticalAlignment = -4108

$ws.Range("D3").Value = 5/24
$ws.Range("D3").NumberFormat = "h:mm"
$ws.Range("D3").HorizontalAlignment = -4108
$ws.Range("D3").VerticalAlignment = -4108

$ws.Range("E3").Value = "1. Loading the display with the menu items`n2. Showing the text in the middle of the screen`n3. Be sure the game screen keeps open when no action performed`n4. Create the menu`n"
$ws.Range("E3").WrapText = $true
$ws.Range("E3").VerticalAlignment = -4160

$ws.Range("F3").Value = "Have a working menu`n"
$ws.Range("F3").WrapText = $true
$ws.Range("F3").VerticalAlignment = -4160

$ws.Rows.Item(3).RowHeight = 93

# -- Row 4 : ID 1, start/stop the game ---------------------------------
$ws.Range("A4").Value = 1
$ws.Range("A4").HorizontalAlignment = -4108
$ws.Range("A4").VerticalAlignment = -4108

$ws.Range("B4").Value = "As a player I want to be able to start and stop the game so that`n I can play the game and end whenever I want"
$ws.Range("B4").WrapText = $true
$ws.Range("B4").VerticalAlignment = -4108     # xlCenter

$ws.Range("C4").Value = "M"
$ws.Range("C4").HorizontalAlignment = -4108
$ws.Range("C4").VerticalAlignment = -4108

$ws.Range("D4").Value = 2/24
$ws.Range("D4").NumberFormat = "h:mm"
$ws.Range("D4").HorizontalAlignment = -4108
$ws.Range("D4").VerticalAlignment = -4108

$ws.Range("E4").Value = "1. Make start button`n2. Make stop button`n3. Close the game when pressed on X"
$ws.Range("E4").WrapText = $true
$ws.Range("E4").VerticalAlignment = -4160

$ws.Range("F4").Value = "Have working start and`nexit buttons"
$ws.Range("F4").WrapText = $true
$ws.Range("F4").VerticalAlignment = -4160

$ws.Rows.Item(4).RowHeight = 63.75

# -- Row 5 : ID 2, view the highscore ----------------------------------
$ws.Range("A5").Value = 2
$ws.Range("A5").HorizontalAlignment = -4108
$ws.Range("A5").VerticalAlignment = -4108

$ws.Range("B5").Value = "As a player I want to view the highscore, so that I can view who has the highest score"
$ws.Range("B5").WrapText = $true
$ws.Range("B5").VerticalAlignment = -4108

$ws.Range("C5").Value = "M"
$ws.Range("C5").HorizontalAlignment = -4108
$ws.Range("C5").VerticalAlignment = -4108

$ws.Range("D5").Value = 2/24
$ws.Range("D5").NumberFormat = "h:mm"
$ws.Range("D5").HorizontalAlignment = -4108
$ws.Range("D5").VerticalAlignment = -4108

$ws.Range("E5").Value = "1. Make a highscore button`n2. Go to highscore page"
$ws.Range("E5").WrapText = $true
$ws.Range("E5").VerticalAlignment = -4160

$ws.Range("F5").Value = "Have a highscore button"
$ws.Range("F5").VerticalAlignment = -4160

$ws.Rows.Item(5).RowHeight = 56.25

# -- Row 6 : ID 6, read the rules --------------------------------------
$ws.Range("A6").Value = 6
$ws.Range("A6").HorizontalAlignment = -4108
$ws.Range("A6").VerticalAlignment = -4108

$ws.Range("B6").Value = "As a player I want to read the rules from the game, so that I know how the game works"
$ws.Range("B6").WrapText = $true
$ws.Range("B6").VerticalAlignment = -4108

$ws.Range("C6").Value = "M"
$ws.Range("C6").HorizontalAlignment = -4108
$ws.Range("C6").VerticalAlignment = -4108

$ws.Range("D6").Value = 1.5/24
$ws.Range("D6").NumberFormat = "h:mm"
$ws.Range("D6").HorizontalAlignment = -4108
$ws.Range("D6").VerticalAlignment = -4108

$ws.Range("E6").Value = "1. Make rules button`n2. Make a rules page"
$ws.Range("E6").WrapText = $true
$ws.Range("E6").VerticalAlignment = -4160

$ws.Range("F6").Value = "Have a rules button"
$ws.Range("F6").VerticalAlignment = -4160

$ws.Rows.Item(6).RowHeight = 54.75

# ---------------------------------------------------------------------
# 4. Column widths for the new table.
# ---------------------------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 40.666666666666664   # -> stored 41.5
$ws.Columns.Item(4).ColumnWidth = 15.166666666666666   # -> stored 16
$ws.Columns.Item(5).ColumnWidth = 39.166666666666664   # -> stored 40
$ws.Columns.Item(6).ColumnWidth = 22                    # -> stored 22.8333...

# ---------------------------------------------------------------------
# 5. Selection state: Sprint1 was left with E4 selected, Product re-
#    gains the active-tab / active-cell state (B19) it had before.
# ---------------------------------------------------------------------
$ws.Range("E4").Select()
$product.Activate()
$product.Range("B19").Select()
